$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 194 (pushing the former
# rows 194..337 down to 195..338, growing the used range to A1:R338).
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new record's data.
$ws.Range("A194").Value = 10
$ws.Range("B194").Value = 'Vega Modelo de Temuco'
$ws.Range("C194").Value = 'La Araucanía'
$ws.Range("D194").Value = 44762
$ws.Range("E194").Value = 9
$ws.Range("F194").Value = 100112044
$ws.Range("G194").Value = 'Perejil'
$ws.Range("H194").Value = 'Sin especificar'
$ws.Range("I194").Value = 'Primera'
$ws.Range("J194").Value = 55
$ws.Range("K194").Value = 4333
$ws.Range("L194").Value = 4333
$ws.Range("M194").Value = 4333
$ws.Range("N194").Value = '$/docena de atados (3 kilos)'
$ws.Range("O194").Value = 'Región Metropolitana'
$ws.Range("P194").Value = 1444
$ws.Range("Q194").Value = 3
$ws.Range("R194").Value = 'Hortaliza'
